$wb = $excel.ActiveWorkbook

# New snapshot row (row 93) appended to each of the 4 sheets, one day
# after the last existing row (row 92). Column A is a datetime serial
# formatted like the existing rows; columns B-E are hex-ish text strings;
# columns F-I are plain numbers.

$rowNum = 93
$dateSerial = 45879.49484953703
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$sheetsData = @(
    @{
        Name = "FE_LFT_#1"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x08"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 264
        I = 15
    },
    @{
        Name = "FE_LFT_#2"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x18"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 280
        I = 14
    },
    @{
        Name = "FE_PLT_#1"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x5D"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 93
        I = 3
    },
    @{
        Name = "FE_PLT_#2"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5B"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 91
        I = 3
    }
)

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Name)

    $ws.Cells.Item($rowNum, 1).Value = $dateSerial
    $ws.Cells.Item($rowNum, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($rowNum, 2).Value = $sd.B
    $ws.Cells.Item($rowNum, 3).Value = $sd.C
    $ws.Cells.Item($rowNum, 4).Value = $sd.D
    $ws.Cells.Item($rowNum, 5).Value = $sd.E

    $ws.Cells.Item($rowNum, 6).Value = $sd.F
    $ws.Cells.Item($rowNum, 7).Value = $sd.G
    $ws.Cells.Item($rowNum, 8).Value = $sd.H
    $ws.Cells.Item($rowNum, 9).Value = $sd.I
}
